# Updated cryptos list — applies Price (D) and Volume(1h) (E) edits
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): text values that must stay text (some look numeric),
# so force text format, assign, then restore the default 'Normal' style
# so no stray number-format style sticks to the cell.
$priceUpdates = @{
    'D2' = '27.029.84'
    'D3' = '1.562.56'
    'D5' = '208.47'
    'D11' = '0.0857'
    'D12' = '1.785.09'
    'D13' = '1.573.54'
    'D14' = '3.74'
    'D16' = '27.025.83'
    'D17' = '61.86'
    'D19' = '215.69'
    'D20' = '7.39'
    'D22' = '4.15'
    'D23' = '9.22'
    'D25' = '153.45'
    'D26' = '6.60'
    'D27' = '15.04'
    'D31' = '1.13'
    'D34' = '1.422.34'
    'D39' = '0.531'
    'D45' = '64.74'
    'D47' = '1.700.08'
    'D48' = '86.59'
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = '@'
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = 'Normal'
}

# Volume(1h) column (E): padded percent strings, already safe as text
$volumeUpdates = @{
    'E2' = '  +0.44%  '
    'E3' = '  +0.81%  '
    'E4' = '  +0.45%  '
    'E5' = '  +0.89%  '
    'E6' = '  +0.76%  '
    'E7' = '  +0.39%  '
    'E8' = '  -0.15%  '
    'E9' = '  +1.52%  '
    'E10' = '  +1.81%  '
    'E11' = '  +0.12%  '
    'E12' = '  +0.80%  '
    'E13' = '  +1.50%  '
    'E15' = '  +0.19%  '
    'E16' = '  +0.47%  '
    'E17' = '  +0.37%  '
    'E18' = '  +1.24%  '
    'E19' = '  -0.81%  '
    'E20' = '  +1.91%  '
    'E21' = '  +0.39%  '
    'E22' = '  +2.35%  '
    'E23' = '  +0.31%  '
    'E24' = '  -0.32%  '
    'E25' = '  -0.51%  '
    'E26' = '  -0.16%  '
    'E27' = '  +0.77%  '
    'E28' = '  +1.43%  '
    'E29' = '  +0.26%  '
    'E30' = '  +1.45%  '
    'E31' = '  +4.12%  '
    'E32' = '  +0.19%  '
    'E33' = '  +3.73%  '
    'E34' = '  +0.46%  '
    'E35' = '  +1.94%  '
    'E36' = '  +9.88%  '
    'E37' = '  +2.60%  '
    'E38' = '  +1.04%  '
    'E39' = '  +1.78%  '
    'E41' = '  +0.31%  '
    'E42' = '  +0.41%  '
    'E43' = '  +0.90%  '
    'E44' = '  -0.11%  '
    'E45' = '  +0.74%  '
    'E46' = '  -0.96%  '
    'E47' = '  +0.95%  '
    'E48' = '  -1.11%  '
    'E49' = '  +3.18%  '
    'E50' = '  -0.37%  '
    'E51' = '  +0.50%  '
}
foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
